# Update LDLC price-tracking workbook: add a new price-snapshot column.
#
# The sheet keeps a running history of prices, one column per scrape
# timestamp, followed by two fixed trailing columns "nom" and
# "url_produit". A new scrape ("2026-02-05 19:31:36") needs to be
# inserted right before the "nom" column (currently GI), pushing
# "nom" / "url_produit" one column to the right (GI->GJ, GJ->GK) and
# carrying the latest known price forward into the freshly inserted
# column for every product row that still has a tracked price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column in front of the "nom" column (GI). This
# shifts the existing "nom"/"url_produit" columns (GI/GJ) one slot to
# the right (-> GJ/GK) and extends the sheet dimension accordingly.
$ws.Range("GI1").EntireColumn.Insert()

# Header for the newly inserted snapshot column.
$ws.Range("GI1").Value = "2026-02-05 19:31:36"

# For every product row that still had a numeric price in the previous
# latest snapshot (column GH, rows 2-80), carry that same price value
# forward into the newly inserted GI column. Rows 81-210 correspond to
# products whose price tracking already stopped (blank snapshot
# cells), so they are left untouched / blank.
For ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, 190).Copy($ws.Cells.Item($r, 191))
}
